$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист3")

# A1:D121 grid data (anode grid positions + hit counts)
$arr = New-Object 'object[,]' 121,4
$arr[0,0] = 0; $arr[0,1] = -50; $arr[0,2] = -50; $arr[0,3] = 0
$arr[1,0] = 1; $arr[1,1] = -40; $arr[1,2] = -50; $arr[1,3] = 0
$arr[2,0] = 2; $arr[2,1] = -30; $arr[2,2] = -50; $arr[2,3] = 0
$arr[3,0] = 3; $arr[3,1] = -20; $arr[3,2] = -50; $arr[3,3] = 0
$arr[4,0] = 4; $arr[4,1] = -10; $arr[4,2] = -50; $arr[4,3] = 0
$arr[5,0] = 5; $arr[5,1] = 0; $arr[5,2] = -50; $arr[5,3] = 0
$arr[6,0] = 6; $arr[6,1] = 10; $arr[6,2] = -50; $arr[6,3] = 0
$arr[7,0] = 7; $arr[7,1] = 20; $arr[7,2] = -50; $arr[7,3] = 0
$arr[8,0] = 8; $arr[8,1] = 30; $arr[8,2] = -50; $arr[8,3] = 0
$arr[9,0] = 9; $arr[9,1] = 40; $arr[9,2] = -50; $arr[9,3] = 0
$arr[10,0] = 10; $arr[10,1] = 50; $arr[10,2] = -50; $arr[10,3] = 0
$arr[11,0] = 11; $arr[11,1] = -50; $arr[11,2] = -40; $arr[11,3] = 0
$arr[12,0] = 12; $arr[12,1] = -40; $arr[12,2] = -40; $arr[12,3] = 0
$arr[13,0] = 13; $arr[13,1] = -30; $arr[13,2] = -40; $arr[13,3] = 0
$arr[14,0] = 14; $arr[14,1] = -20; $arr[14,2] = -40; $arr[14,3] = 0
$arr[15,0] = 15; $arr[15,1] = -10; $arr[15,2] = -40; $arr[15,3] = 0
$arr[16,0] = 16; $arr[16,1] = 0; $arr[16,2] = -40; $arr[16,3] = 0
$arr[17,0] = 17; $arr[17,1] = 10; $arr[17,2] = -40; $arr[17,3] = 0
$arr[18,0] = 18; $arr[18,1] = 20; $arr[18,2] = -40; $arr[18,3] = 0
$arr[19,0] = 19; $arr[19,1] = 30; $arr[19,2] = -40; $arr[19,3] = 0
$arr[20,0] = 20; $arr[20,1] = 40; $arr[20,2] = -40; $arr[20,3] = 0
$arr[21,0] = 21; $arr[21,1] = 50; $arr[21,2] = -40; $arr[21,3] = 0
$arr[22,0] = 22; $arr[22,1] = -50; $arr[22,2] = -30; $arr[22,3] = 0
$arr[23,0] = 23; $arr[23,1] = -40; $arr[23,2] = -30; $arr[23,3] = 0
$arr[24,0] = 24; $arr[24,1] = -30; $arr[24,2] = -30; $arr[24,3] = 0
$arr[25,0] = 25; $arr[25,1] = -20; $arr[25,2] = -30; $arr[25,3] = 0
$arr[26,0] = 26; $arr[26,1] = -10; $arr[26,2] = -30; $arr[26,3] = 0
$arr[27,0] = 27; $arr[27,1] = 0; $arr[27,2] = -30; $arr[27,3] = 1
$arr[28,0] = 28; $arr[28,1] = 10; $arr[28,2] = -30; $arr[28,3] = 0
$arr[29,0] = 29; $arr[29,1] = 20; $arr[29,2] = -30; $arr[29,3] = 0
$arr[30,0] = 30; $arr[30,1] = 30; $arr[30,2] = -30; $arr[30,3] = 0
$arr[31,0] = 31; $arr[31,1] = 40; $arr[31,2] = -30; $arr[31,3] = 0
$arr[32,0] = 32; $arr[32,1] = 50; $arr[32,2] = -30; $arr[32,3] = 0
$arr[33,0] = 33; $arr[33,1] = -50; $arr[33,2] = -20; $arr[33,3] = 0
$arr[34,0] = 34; $arr[34,1] = -40; $arr[34,2] = -20; $arr[34,3] = 0
$arr[35,0] = 35; $arr[35,1] = -30; $arr[35,2] = -20; $arr[35,3] = 0
$arr[36,0] = 36; $arr[36,1] = -20; $arr[36,2] = -20; $arr[36,3] = 0
$arr[37,0] = 37; $arr[37,1] = -10; $arr[37,2] = -20; $arr[37,3] = 0
$arr[38,0] = 38; $arr[38,1] = 0; $arr[38,2] = -20; $arr[38,3] = 0
$arr[39,0] = 39; $arr[39,1] = 10; $arr[39,2] = -20; $arr[39,3] = 0
$arr[40,0] = 40; $arr[40,1] = 20; $arr[40,2] = -20; $arr[40,3] = 0
$arr[41,0] = 41; $arr[41,1] = 30; $arr[41,2] = -20; $arr[41,3] = 0
$arr[42,0] = 42; $arr[42,1] = 40; $arr[42,2] = -20; $arr[42,3] = 0
$arr[43,0] = 43; $arr[43,1] = 50; $arr[43,2] = -20; $arr[43,3] = 0
$arr[44,0] = 44; $arr[44,1] = -50; $arr[44,2] = -10; $arr[44,3] = 0
$arr[45,0] = 45; $arr[45,1] = -40; $arr[45,2] = -10; $arr[45,3] = 0
$arr[46,0] = 46; $arr[46,1] = -30; $arr[46,2] = -10; $arr[46,3] = 0
$arr[47,0] = 47; $arr[47,1] = -20; $arr[47,2] = -10; $arr[47,3] = 0
$arr[48,0] = 48; $arr[48,1] = -10; $arr[48,2] = -10; $arr[48,3] = 1
$arr[49,0] = 49; $arr[49,1] = 0; $arr[49,2] = -10; $arr[49,3] = 0
$arr[50,0] = 50; $arr[50,1] = 10; $arr[50,2] = -10; $arr[50,3] = 0
$arr[51,0] = 51; $arr[51,1] = 20; $arr[51,2] = -10; $arr[51,3] = 0
$arr[52,0] = 52; $arr[52,1] = 30; $arr[52,2] = -10; $arr[52,3] = 0
$arr[53,0] = 53; $arr[53,1] = 40; $arr[53,2] = -10; $arr[53,3] = 0
$arr[54,0] = 54; $arr[54,1] = 50; $arr[54,2] = -10; $arr[54,3] = 0
$arr[55,0] = 55; $arr[55,1] = -50; $arr[55,2] = 0; $arr[55,3] = 0
$arr[56,0] = 56; $arr[56,1] = -40; $arr[56,2] = 0; $arr[56,3] = 0
$arr[57,0] = 57; $arr[57,1] = -30; $arr[57,2] = 0; $arr[57,3] = 0
$arr[58,0] = 58; $arr[58,1] = -20; $arr[58,2] = 0; $arr[58,3] = 0
$arr[59,0] = 59; $arr[59,1] = -10; $arr[59,2] = 0; $arr[59,3] = 3
$arr[60,0] = 60; $arr[60,1] = 0; $arr[60,2] = 0; $arr[60,3] = 10
$arr[61,0] = 61; $arr[61,1] = 10; $arr[61,2] = 0; $arr[61,3] = 1
$arr[62,0] = 62; $arr[62,1] = 20; $arr[62,2] = 0; $arr[62,3] = 0
$arr[63,0] = 63; $arr[63,1] = 30; $arr[63,2] = 0; $arr[63,3] = 0
$arr[64,0] = 64; $arr[64,1] = 40; $arr[64,2] = 0; $arr[64,3] = 0
$arr[65,0] = 65; $arr[65,1] = 50; $arr[65,2] = 0; $arr[65,3] = 0
$arr[66,0] = 66; $arr[66,1] = -50; $arr[66,2] = 10; $arr[66,3] = 0
$arr[67,0] = 67; $arr[67,1] = -40; $arr[67,2] = 10; $arr[67,3] = 0
$arr[68,0] = 68; $arr[68,1] = -30; $arr[68,2] = 10; $arr[68,3] = 0
$arr[69,0] = 69; $arr[69,1] = -20; $arr[69,2] = 10; $arr[69,3] = 0
$arr[70,0] = 70; $arr[70,1] = -10; $arr[70,2] = 10; $arr[70,3] = 1
$arr[71,0] = 71; $arr[71,1] = 0; $arr[71,2] = 10; $arr[71,3] = 1
$arr[72,0] = 72; $arr[72,1] = 10; $arr[72,2] = 10; $arr[72,3] = 0
$arr[73,0] = 73; $arr[73,1] = 20; $arr[73,2] = 10; $arr[73,3] = 0
$arr[74,0] = 74; $arr[74,1] = 30; $arr[74,2] = 10; $arr[74,3] = 0
$arr[75,0] = 75; $arr[75,1] = 40; $arr[75,2] = 10; $arr[75,3] = 0
$arr[76,0] = 76; $arr[76,1] = 50; $arr[76,2] = 10; $arr[76,3] = 0
$arr[77,0] = 77; $arr[77,1] = -50; $arr[77,2] = 20; $arr[77,3] = 0
$arr[78,0] = 78; $arr[78,1] = -40; $arr[78,2] = 20; $arr[78,3] = 0
$arr[79,0] = 79; $arr[79,1] = -30; $arr[79,2] = 20; $arr[79,3] = 0
$arr[80,0] = 80; $arr[80,1] = -20; $arr[80,2] = 20; $arr[80,3] = 1
$arr[81,0] = 81; $arr[81,1] = -10; $arr[81,2] = 20; $arr[81,3] = 0
$arr[82,0] = 82; $arr[82,1] = 0; $arr[82,2] = 20; $arr[82,3] = 0
$arr[83,0] = 83; $arr[83,1] = 10; $arr[83,2] = 20; $arr[83,3] = 0
$arr[84,0] = 84; $arr[84,1] = 20; $arr[84,2] = 20; $arr[84,3] = 0
$arr[85,0] = 85; $arr[85,1] = 30; $arr[85,2] = 20; $arr[85,3] = 0
$arr[86,0] = 86; $arr[86,1] = 40; $arr[86,2] = 20; $arr[86,3] = 0
$arr[87,0] = 87; $arr[87,1] = 50; $arr[87,2] = 20; $arr[87,3] = 0
$arr[88,0] = 88; $arr[88,1] = -50; $arr[88,2] = 30; $arr[88,3] = 0
$arr[89,0] = 89; $arr[89,1] = -40; $arr[89,2] = 30; $arr[89,3] = 0
$arr[90,0] = 90; $arr[90,1] = -30; $arr[90,2] = 30; $arr[90,3] = 0
$arr[91,0] = 91; $arr[91,1] = -20; $arr[91,2] = 30; $arr[91,3] = 0
$arr[92,0] = 92; $arr[92,1] = -10; $arr[92,2] = 30; $arr[92,3] = 0
$arr[93,0] = 93; $arr[93,1] = 0; $arr[93,2] = 30; $arr[93,3] = 0
$arr[94,0] = 94; $arr[94,1] = 10; $arr[94,2] = 30; $arr[94,3] = 0
$arr[95,0] = 95; $arr[95,1] = 20; $arr[95,2] = 30; $arr[95,3] = 0
$arr[96,0] = 96; $arr[96,1] = 30; $arr[96,2] = 30; $arr[96,3] = 0
$arr[97,0] = 97; $arr[97,1] = 40; $arr[97,2] = 30; $arr[97,3] = 0
$arr[98,0] = 98; $arr[98,1] = 50; $arr[98,2] = 30; $arr[98,3] = 0
$arr[99,0] = 99; $arr[99,1] = -50; $arr[99,2] = 40; $arr[99,3] = 0
$arr[100,0] = 100; $arr[100,1] = -40; $arr[100,2] = 40; $arr[100,3] = 0
$arr[101,0] = 101; $arr[101,1] = -30; $arr[101,2] = 40; $arr[101,3] = 0
$arr[102,0] = 102; $arr[102,1] = -20; $arr[102,2] = 40; $arr[102,3] = 0
$arr[103,0] = 103; $arr[103,1] = -10; $arr[103,2] = 40; $arr[103,3] = 0
$arr[104,0] = 104; $arr[104,1] = 0; $arr[104,2] = 40; $arr[104,3] = 0
$arr[105,0] = 105; $arr[105,1] = 10; $arr[105,2] = 40; $arr[105,3] = 0
$arr[106,0] = 106; $arr[106,1] = 20; $arr[106,2] = 40; $arr[106,3] = 0
$arr[107,0] = 107; $arr[107,1] = 30; $arr[107,2] = 40; $arr[107,3] = 0
$arr[108,0] = 108; $arr[108,1] = 40; $arr[108,2] = 40; $arr[108,3] = 0
$arr[109,0] = 109; $arr[109,1] = 50; $arr[109,2] = 40; $arr[109,3] = 0
$arr[110,0] = 110; $arr[110,1] = -50; $arr[110,2] = 50; $arr[110,3] = 0
$arr[111,0] = 111; $arr[111,1] = -40; $arr[111,2] = 50; $arr[111,3] = 0
$arr[112,0] = 112; $arr[112,1] = -30; $arr[112,2] = 50; $arr[112,3] = 0
$arr[113,0] = 113; $arr[113,1] = -20; $arr[113,2] = 50; $arr[113,3] = 0
$arr[114,0] = 114; $arr[114,1] = -10; $arr[114,2] = 50; $arr[114,3] = 0
$arr[115,0] = 115; $arr[115,1] = 0; $arr[115,2] = 50; $arr[115,3] = 0
$arr[116,0] = 116; $arr[116,1] = 10; $arr[116,2] = 50; $arr[116,3] = 0
$arr[117,0] = 117; $arr[117,1] = 20; $arr[117,2] = 50; $arr[117,3] = 0
$arr[118,0] = 118; $arr[118,1] = 30; $arr[118,2] = 50; $arr[118,3] = 0
$arr[119,0] = 119; $arr[119,1] = 40; $arr[119,2] = 50; $arr[119,3] = 0
$arr[120,0] = 120; $arr[120,1] = 50; $arr[120,2] = 50; $arr[120,3] = 0
$ws.Range("A1:D121").Value = $arr

# H1:M5 summary statistics block
$ws.Range("H1").Value = "N_emit"
$ws.Range("I1").Value = "N_reg_tot"
$ws.Range("J1").Value = "x_avr"
$ws.Range("K1").Value = "y_avr"
$ws.Range("L1").Value = "x_recon"
$ws.Range("M1").Value = "y_recon"
$ws.Range("H2").Value = 1600
$ws.Range("I2").Value = 19
$ws.Range("J2").Value = 0.45
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = -28.8276
$ws.Range("M2").Value = -5.24138
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 20.7273
$ws.Range("L4").Value = -3.30435
$ws.Range("M4").Value = 6.6087
$ws.Range("L5").Value = 10.1333
$ws.Range("M5").Value = -5.06667

# Column I best-fit width
$ws.Columns.Item(9).AutoFit()

# Switch active sheet to Лист3 and set selection
[void]$ws.Select()
[void]$ws.Range("I6").Select()
